$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.029.74"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3
$ws.Range("D3").Value = "1.651.66"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'218.36"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  +0.50%  "

# Row 7
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  +1.80%  "

# Row 9
$ws.Range("D9").Value = "'0.0623"
$ws.Range("E9").Value = "  +0.14%  "

# Row 10
$ws.Range("D10").Value = "'19.80"
$ws.Range("E10").Value = "  +4.11%  "

# Row 11
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").Value = "1.886.52"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("D13").Value = "1.655.40"
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").Value = "'4.16"
$ws.Range("E14").Value = "  +0.59%  "

# Row 15
$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").Value = "'66.75"
$ws.Range("E16").Value = "  +3.81%  "

# Row 17
$ws.Range("D17").Value = "27.098.24"
$ws.Range("E17").Value = "  +1.50%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").Value = "'220.50"
$ws.Range("E19").Value = "  +4.90%  "

# Row 20
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  +8.34%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.42"
$ws.Range("E22").Value = "  +2.07%  "

# Row 23
$ws.Range("D23").Value = "'2.40"
$ws.Range("E23").Value = "  +3.28%  "

# Row 24
$ws.Range("D24").Value = "'9.19"
$ws.Range("E24").Value = "  -0.39%  "

# Row 25
$ws.Range("D25").Value = "'146.75"
$ws.Range("E25").Value = "  +0.85%  "

# Row 26
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27
$ws.Range("D27").Value = "'7.42"
$ws.Range("E27").Value = "  +5.13%  "

# Row 28
$ws.Range("E28").Value = "  +1.40%  "

# Row 29
$ws.Range("D29").Value = "'15.95"
$ws.Range("E29").Value = "  +2.68%  "

# Row 30
$ws.Range("D30").Value = "'0.0513"
$ws.Range("E30").Value = "  +1.77%  "

# Row 31
$ws.Range("D31").Value = "'1.19"
$ws.Range("E31").Value = "  +0.86%  "

# Row 32
$ws.Range("D32").Value = "'3.40"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  +0.80%  "

# Row 34
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  +2.38%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.254.73"
$ws.Range("E36").Value = "  -1.43%  "

# Row 37
$ws.Range("D37").Value = "'0.0177"
$ws.Range("E37").Value = "  +1.38%  "

# Row 38
$ws.Range("D38").Value = "'0.534"
$ws.Range("E38").Value = "  +0.78%  "

# Row 39
$ws.Range("D39").Value = "'0.833"
$ws.Range("E39").Value = "  +3.14%  "

# Row 40
$ws.Range("D40").Value = "'1.01"
$ws.Range("E40").Value = "  -0.11%  "

# Row 41
$ws.Range("D41").Value = "'0.815"
$ws.Range("E41").Value = "  +1.67%  "

# Row 42
$ws.Range("D42").Value = "'5.36"
$ws.Range("E42").Value = "  +2.09%  "

# Row 43
$ws.Range("D43").Value = "1.797.76"
$ws.Range("E43").Value = "  +1.33%  "

# Row 44
$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  -4.17%  "

# Row 45
$ws.Range("D45").Value = "'61.46"
$ws.Range("E45").Value = "  +1.74%  "

# Row 46
$ws.Range("D46").Value = "'91.52"
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = "  +2.05%  "

# Row 48
$ws.Range("E48").Value = "  -0.69%  "

# Row 49
$ws.Range("D49").Value = "'0.0978"
$ws.Range("E49").Value = "  +1.86%  "

# Row 50
$ws.Range("D50").Value = "'7.64"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0971"
$ws.Range("E51").Value = "  -4.86%  "
